# Update "想去人数" (interest count) figures for a few entries that changed
# between scrapes, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1658
    $ws.Range("F5").Value = 6215
    $ws.Range("F6").Value = 45
}
